$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41, shifting Gomez-Cadenas and everything
# below it down by one (matches the diff: dimension A1:H110 -> A1:H111).
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with the new author entry (Garcia-Barrena, UPV Valencia).
$ws.Cells.Item(41, 1).Value = 'Garc\''ia-Barrena'
$ws.Cells.Item(41, 2).Value = 'J.'
$ws.Cells.Item(41, 5).Value = 'Instituto de Instrumentaci\''on para Imagen Molecular (I3M), Centro Mixto CSIC - Universitat Polit\`ecnica de Val\`encia, Camino de Vera s/n '
$ws.Cells.Item(41, 6).Value = ' Valencia, E-46022, Spain'

# Update the view's selection to match the author's edit position (cell A41,
# where the new row was inserted). The engine doesn't model a Window object
# for the topLeftCell scroll position, so only the selection can be set here.
$ws.Range("A41").Select()
